$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 493, pushing the existing row 493 (and all
# rows below it) down by one.
$ws.Rows(493).Insert()

# Populate the newly inserted row 493 with the new data record (same
# categorical data as the original row that is now at 494, with updated
# date / price figures).
$ws.Range("A493").Value = 9
$ws.Range("B493").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C493").Value = "Metropolitana"
$ws.Range("D493").Value = 44946
$ws.Range("E493").Value = 13
$ws.Range("F493").Value = 100112012
$ws.Range("G493").Value = "Espinaca"
$ws.Range("H493").Value = "Sin especificar"
$ws.Range("I493").Value = "Primera"
$ws.Range("J493").Value = 160
$ws.Range("K493").Value = 6000
$ws.Range("L493").Value = 8000
$ws.Range("M493").Value = 7000
$ws.Range("N493").Value = "$/cuna 10 kilos"
$ws.Range("O493").Value = "Provincia de Chacabuco"
$ws.Range("P493").Value = 700
$ws.Range("Q493").Value = 10
$ws.Range("R493").Value = "Hortaliza"
